$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.771.80"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.19"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.60%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("E7").Value = "  +0.85%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3930"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07929"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9777"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.91%  "

$ws.Range("E11").Value = "  +2.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.887.82"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.741"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.013"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06953"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.23%  "

$ws.Range("E16").Value = "  +2.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("E18").Value = "  +1.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.97"
$ws.Range("D19").ClearFormats()

$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.772.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.355"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.08"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.121"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.116.64"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.60"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.07%  "

$ws.Range("E27").Value = "  +1.19%  "

$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.998"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.96"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.86%  "

$ws.Range("E31").Value = "  +1.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9383"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.306"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.355"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.349"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05916"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.157"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.901"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5713"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1798"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("E42").Value = "  +0.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07315"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.85"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.69%  "

$ws.Range("E45").Value = "  +1.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.150"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.109"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.37%  "

$ws.Range("E48").Value = "  +1.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.00"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.369"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.51%  "

$ws.Range("E51").Value = "  +0.42%  "
